$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last refreshed" shared-string cell
$ws.Range("A300").Value = "Ostatnia aktualizacja: 25-07-2023, 10:39"

# Fill in the newly-reported "S" column (latest year) data points,
# copying number formatting from the neighbouring "R" column cell
function Set-SValue($row, $value) {
    $src = $ws.Range("R" + $row)
    $dst = $ws.Range("S" + $row)
    $src.Copy() | Out-Null
    $dst.PasteSpecial(-4122) | Out-Null
    $dst.Value = $value
}

Set-SValue 9 412.4
Set-SValue 10 469.6
Set-SValue 11 328.2
Set-SValue 68 36.1
Set-SValue 82 7.6
Set-SValue 83 11.3
Set-SValue 84 9.4
Set-SValue 85 6.4
Set-SValue 86 9
Set-SValue 87 6.6
Set-SValue 88 2.9
Set-SValue 89 3.3
Set-SValue 129 71.900000000000006
Set-SValue 130 77.3
Set-SValue 131 65.099999999999994
Set-SValue 132 859996.9
Set-SValue 146 56.3
Set-SValue 147 64
Set-SValue 148 49.3
Set-SValue 149 27.8
Set-SValue 150 71.3
Set-SValue 151 77.8
Set-SValue 152 51
Set-SValue 153 76.7
Set-SValue 154 35.1
Set-SValue 155 56.4
Set-SValue 156 12.3
Set-SValue 157 56.3
Set-SValue 158 56.3
Set-SValue 159 80.3
Set-SValue 160 59.9
Set-SValue 161 52.5
Set-SValue 162 53.2
Set-SValue 163 16.3
Set-SValue 167 58
Set-SValue 168 50.8
Set-SValue 169 65.8
Set-SValue 170 57.9
Set-SValue 171 58.1
Set-SValue 172 31.1
Set-SValue 173 80.2
Set-SValue 174 35.799999999999997
Set-SValue 175 20.100000000000001
Set-SValue 176 22.1
Set-SValue 177 21.3
Set-SValue 178 22.8
Set-SValue 179 21
Set-SValue 180 23.6
Set-SValue 181 8
Set-SValue 182 8.6
Set-SValue 183 7.5
Set-SValue 184 1939
Set-SValue 185 2234
Set-SValue 186 -295
Set-SValue 207 30.4
Set-SValue 225 6.6
Set-SValue 246 859996.9
Set-SValue 258 0.19800000000000001
Set-SValue 259 31
Set-SValue 289 4
Set-SValue 292 0.53
Set-SValue 294 9.59

$excel.CutCopyMode = 0

# Selection moved to A4 and print scale bumped 63% -> 65%
$ws.Range("A4").Select() | Out-Null
$ws.PageSetup.Zoom = 65
